$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix row 343: C/D/E/F change from 41569200000 to 37573500000 ---
$ws.Range("C343").Value = 37573500000
$ws.Range("D343").Value = 37573500000
$ws.Range("E343").Value = 37573500000
$ws.Range("F343").Value = 37573500000

# --- Append three new data rows (353-355), matching row 343's formatting ---
$ws.Range("A343").Copy()
$ws.Range("A353:A355").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

$ws.Range("A353").Value = 44986.45833333334
$ws.Range("B353").Value = "ECONOMICS:KWM2"
$ws.Range("C353").Value = 39530200000
$ws.Range("D353").Value = 39530200000
$ws.Range("E353").Value = 39530200000
$ws.Range("F353").Value = 39530200000
$ws.Range("G353").Value = 0

$ws.Range("A354").Value = 45017.45833333334
$ws.Range("B354").Value = "ECONOMICS:KWM2"
$ws.Range("C354").Value = 39461400000
$ws.Range("D354").Value = 39461400000
$ws.Range("E354").Value = 39461400000
$ws.Range("F354").Value = 39461400000
$ws.Range("G354").Value = 0

$ws.Range("A355").Value = 45047.41666666666
$ws.Range("B355").Value = "ECONOMICS:KWM2"
$ws.Range("C355").Value = 39655200000
$ws.Range("D355").Value = 39655200000
$ws.Range("E355").Value = 39655200000
$ws.Range("F355").Value = 39655200000
$ws.Range("G355").Value = 0
